$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 51666.668
$ws.Range("I61").Value = 51666.668
$ws.Range("K61").Value = 51666.668
$ws.Range("M61").Value = -51454.668

$ws.Range("H63").Value = 3450
$ws.Range("J63").Value = 4000
$ws.Range("L63").Value = 4000
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 3450
$ws.Range("J66").Value = 4000
$ws.Range("L66").Value = 20000
$ws.Range("N66").Value = -26864

$ws.Range("H97").Value = 425
$ws.Range("I97").Value = 425
$ws.Range("K97").Value = 425
$ws.Range("M97").Value = 71

$ws.Range("H136").Value = 51666.668
$ws.Range("I136").Value = 51666.668
$ws.Range("K136").Value = 155000.004
$ws.Range("M136").Value = -152450.004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 241.16667
$ws.Range("I10").Value = 284.66666
$ws.Range("J10").Value = 197.66667
$ws.Range("K10").Value = 284.66666
$ws.Range("L10").Value = 197.66667
$ws.Range("M10").Value = -144.66666
$ws.Range("N10").Value = -477.66667

$ws.Range("H11").Value = 7500
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 11000
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 11000
$ws.Range("M11").Value = -360
$ws.Range("N11").Value = -11280

$ws.Range("H12").Value = 1868.3334
$ws.Range("J12").Value = 2552.5
$ws.Range("L12").Value = 2552.5
$ws.Range("N12").Value = -2888.5

$ws.Range("H14").Value = 6000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H22").Value = 1220

$ws.Range("H99").Value = 990
$ws.Range("I99").Value = 970
$ws.Range("K99").Value = 970
$ws.Range("M99").Value = 528

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 9166.666999999999
$ws.Range("I17").Value = 4666.6665
$ws.Range("J17").Value = 13666.667
$ws.Range("K17").Value = 4666.6665
$ws.Range("L17").Value = 13666.667
$ws.Range("M17").Value = -4492.6665
$ws.Range("N17").Value = -14014.667

$ws.Range("H22").Value = 718.7692
$ws.Range("I22").Value = 686.2727
$ws.Range("K22").Value = 686.2727
$ws.Range("M22").Value = -336.2727

$ws.Range("H28").Value = 51175.89
$ws.Range("J28").Value = 51175.89
$ws.Range("L28").Value = 51175.89
$ws.Range("N28").Value = -51665.89

$ws.Range("H31").Value = 3060.8948
$ws.Range("I31").Value = 3073.5
$ws.Range("K31").Value = 3073.5
$ws.Range("M31").Value = -2778.5

$ws.Range("H34").Value = 3060.8948
$ws.Range("I34").Value = 3073.5
$ws.Range("K34").Value = 3073.5
$ws.Range("M34").Value = -2871.5

$ws.Range("H94").Value = 8250
$ws.Range("I94").Value = 8250
$ws.Range("K94").Value = 8250
$ws.Range("M94").Value = -7799

$ws.Range("H105").Value = 1225
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 468.5
$ws.Range("J23").Value = 594.6667
$ws.Range("L23").Value = 1784.0001
$ws.Range("N23").Value = -2254.0001

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H86").Value = 561.75
$ws.Range("J86").Value = 561.75
$ws.Range("L86").Value = 1685.25
$ws.Range("N86").Value = -4057.25

$ws.Range("H89").Value = 561.75
$ws.Range("J89").Value = 561.75
$ws.Range("L89").Value = 5055.75
$ws.Range("N89").Value = -16911.75

$ws.Range("H98").Value = 697
$ws.Range("J98").Value = 760.6667
$ws.Range("L98").Value = 2282.0001
$ws.Range("N98").Value = -5278.0001

$ws.Range("H119").Value = 3579.6
$ws.Range("I119").Value = 3579.6
$ws.Range("K119").Value = 10738.8
$ws.Range("M119").Value = -5900.799999999999

$ws.Range("H139").Value = 2521.8333
$ws.Range("I139").Value = 2521.8333
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 7565.499899999999
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -2425.499899999999
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1139.2858
$ws.Range("I13").Value = 250
$ws.Range("J13").Value = 1495
$ws.Range("K13").Value = 250
$ws.Range("L13").Value = 1495
$ws.Range("M13").Value = -111
$ws.Range("N13").Value = -1773

$ws.Range("H97").Value = 2126.4
$ws.Range("I97").Value = 1840
$ws.Range("J97").Value = 2556
$ws.Range("K97").Value = 1840
$ws.Range("L97").Value = 2556
$ws.Range("M97").Value = -1344
$ws.Range("N97").Value = -3548

$ws.Range("H118").Value = 189655
$ws.Range("J118").Value = 189655
$ws.Range("L118").Value = 189655
$ws.Range("N118").Value = -192969

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -888

$ws.Range("H15").Value = 1000
$ws.Range("I15").Value = 1000
$ws.Range("K15").Value = 1000
$ws.Range("M15").Value = -830

$ws.Range("H22").Value = 820.75
$ws.Range("I22").Value = 742
$ws.Range("J22").Value = 899.5
$ws.Range("K22").Value = 742
$ws.Range("L22").Value = 899.5
$ws.Range("M22").Value = -447
$ws.Range("N22").Value = -1489.5

$ws.Range("H27").Value = 820.75
$ws.Range("I27").Value = 742
$ws.Range("J27").Value = 899.5
$ws.Range("K27").Value = 742
$ws.Range("L27").Value = 899.5
$ws.Range("M27").Value = -635
$ws.Range("N27").Value = -1113.5

$ws.Range("H40").Value = 2855.4285
$ws.Range("I40").Value = 1999.6
$ws.Range("J40").Value = 4995
$ws.Range("K40").Value = 1999.6
$ws.Range("L40").Value = 4995
$ws.Range("M40").Value = -1863.6
$ws.Range("N40").Value = -5267

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H68").Value = 1000
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 1000
$ws.Range("M68").Value = -251

$ws.Range("H71").Value = 1000
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 5000
$ws.Range("M71").Value = -1256

$ws.Range("H100").Value = 1036.25
$ws.Range("I100").Value = 1083.3334
$ws.Range("J100").Value = 895
$ws.Range("K100").Value = 1083.3334
$ws.Range("L100").Value = 895
$ws.Range("M100").Value = -542.3334
$ws.Range("N100").Value = -1977

$ws.Range("H104").Value = 58661.332
$ws.Range("J104").Value = 58661.332
$ws.Range("L104").Value = 58661.332
$ws.Range("N104").Value = -65649.33199999999

$ws.Range("H122").Value = 17994.6
$ws.Range("J122").Value = 19992.5
$ws.Range("L122").Value = 59977.5
$ws.Range("N122").Value = -64877.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H4").Value = 14154
$ws.Range("J4").Value = 14154
$ws.Range("L4").Value = 14154
$ws.Range("N4").Value = -14380

$ws.Range("H96").Value = 2420.2856
$ws.Range("I96").Value = 2123.25
$ws.Range("J96").Value = 2816.3333
$ws.Range("K96").Value = 2123.25
$ws.Range("L96").Value = 2816.3333
$ws.Range("M96").Value = -750.25
$ws.Range("N96").Value = -5562.3333

$ws.Range("H122").Value = 17449.666
$ws.Range("I122").Value = 1234.6666
$ws.Range("J122").Value = 33664.668
$ws.Range("K122").Value = 3703.9998
$ws.Range("L122").Value = 100994.004
$ws.Range("M122").Value = -1253.9998
$ws.Range("N122").Value = -105894.004

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
